# Add an "Abbrev" column (C) with short region abbreviations next to the
# existing Code/Name columns on the BaseRegions sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("C1").Value = "Abbrev"
$ws.Range("C1").Font.Bold = $true

# Data rows
$ws.Range("C2").Value = "EU"
$ws.Range("C3").Value = "CN"
$ws.Range("C4").Value = "BR"
$ws.Range("C5").Value = "AU"
$ws.Range("C6").Value = "RoW"

# Move the active selection to the new header cell
$ws.Range("C1").Select() | Out-Null

# Match the print/page setup that Excel records when the sheet is touched
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
